$d = $word.ActiveDocument

# 1. Diseño de interfaz: €2,500 -> €1,500
$d.Content.Find.Execute("Diseño de interfaz: €2,500", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Diseño de interfaz: €1,500", 2)

# 2. Programación y codificación: €6,000 -> €2,500
$d.Content.Find.Execute("Programación y codificación: €6,000", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Programación y codificación: €2,500", 2)

# 3. Remove the whole "Pruebas y control de calidad: €1,500" paragraph
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Pruebas y control de calidad: €1,500*") {
        $p.Range.Delete()
    }
}

# 4. Capacitación y documentación: €2,000 -> €1,000
$d.Content.Find.Execute("Capacitación y documentación: €2,000", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Capacitación y documentación: €1,000", 2)

# 5. Remove the whole "Salarios del equipo del proyecto: €10,000" paragraph
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Salarios del equipo del proyecto: €10,000*") {
        $p.Range.Delete()
    }
}

# 6. Costos Adicionales (imprevistos): €2,000 -> €1,000
$d.Content.Find.Execute("Costos Adicionales (imprevistos): €2,000", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Costos Adicionales (imprevistos): €1,000", 2)

# 7. Insert two new empty paragraphs right after the last table, before the
#    existing trailing empty paragraph (which sits right before the sectPr).
$lastTable = $d.Tables.Item($d.Tables.Count)
$insertRange = $d.Range($lastTable.Range.End, $lastTable.Range.End)
$insertRange.InsertParagraphBefore()
$insertRange = $d.Range($lastTable.Range.End, $lastTable.Range.End)
$insertRange.InsertParagraphBefore()
